$d = $word.ActiveDocument

# 1. Heading "Вывод:" -> "Závěr:"
$d.Content.Find.Execute("Вывод:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Závěr:", 2) | Out-Null

# 2. Replace the old conclusion paragraph text with the new one.
$old = "Серьёзной проблемой для всех беспроводных локальных сетей является безопасность. Она здесь так же важна, как и для любого пользователя сети Интернет. Безопасность является сложным вопросом и требует постоянного внимания. Огромный вред может быть нанесен пользователю из-за того, что он использует случайные хот-споты (hot-spot) или открытые точки доступа WI-FI дома или в офисе и не использует шифрование или VPN (Virtual Private Network - виртуальная частная сеть). Опасно это тем, что пользователь вводит свои личные или профессиональные данные, а сеть при этом не защищена от постороннего вторжения."
$new = "Безопасность беспроводной сети менялась с течением времени, чтобы стать более надёжной, но при этом и более простой с точки зрения её настройки. Но как бы не улучшались протоколы самой главной уязвимость всегда будет оставаться человек. Никто не застрахован от взлома путём фишинга или социальной инженерии, но если вы будете придерживаться выше перечисленных принципов безопасности, то это поможет защитить вашу сеть от взлома."

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new, 2) | Out-Null

# 3. Move the "_GoBack" bookmark from its old spot to sit right after
#    "путём" in the rewritten paragraph (matching the author's edit position).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$marker = "путём"
$markerRange = $d.Content
$found2 = $markerRange.Find.Execute($marker, $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
if ($found2) {
    $pos = $markerRange.End
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

Write-Host "done"
